# Update column G ("K") values on the active sheet.
# The save_data regeneration now uses K (strike count) values computed
# from the freshly calculated std/mean based s_vals instead of the old
# "Strike#" values. Write the newly calculated values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 0
    26 = 3
    27 = 2
    28 = 2
    30 = 1
    31 = 3
    32 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
